# Apply a permutation of the (Fecha, Calidad, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado, Precio $/Kg) tuple across
# rows 2-17 of the active sheet, as described by the commit diff.
# Columns: D=Fecha, L=Calidad, M=Volumen, N=Precio minimo, O=Precio maximo,
#          P=Precio promedio ponderado, S=Precio $/Kg

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: source row -> destination row (value currently in source row
# ends up in destination row after the edit)
$mapping = @{
    2  = 10
    3  = 2
    4  = 13
    5  = 14
    6  = 16
    7  = 17
    8  = 4
    9  = 12
    10 = 15
    11 = 7
    12 = 8
    13 = 3
    14 = 11
    15 = 5
    16 = 6
    17 = 9
}

$columns = @("D", "L", "M", "N", "O", "P", "S")

# Snapshot current values first, since we will overwrite rows while reading.
# Use Value2 (not Value) so dates/numbers round-trip as plain numerics
# rather than being coerced into formatted/display strings.
$snapshot = @{}
foreach ($row in $mapping.Keys) {
    $rowValues = @{}
    foreach ($col in $columns) {
        $rowValues[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $rowValues
}

# Write snapshotted values from each source row into its destination row
foreach ($row in $mapping.Keys) {
    $destRow = $mapping[$row]
    $rowValues = $snapshot[$row]
    foreach ($col in $columns) {
        $ws.Range("$col$destRow").Value2 = $rowValues[$col]
    }
}
